# Generate Report for Handoff
# The localization-status report is regenerated: the file "0bcd2ce7-..."
# moves from "In Translation" to "Ready for handoff" (with new handoff
# timestamps), and the three file rows are re-sorted so that
# "b7d44263-...", "da5e8a2c-..." and "0bcd2ce7-..." appear in that order
# (instead of the original "0bcd2ce7-...", "b7d44263-...", "da5e8a2c-...").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "b7d44263-9138-4e17-8083-9580d31ed741.md"
$ov.Range("B2").Value = "In Translation"
$ov.Range("C2").Value = "In Translation"

$ov.Range("A3").Value = "da5e8a2c-8d37-416c-94e5-65f807b540ad.md"
$ov.Range("B3").Value = "In Translation"
$ov.Range("C3").Value = "In Translation"

$ov.Range("A4").Value = "0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"

$ov.Range("A5").Value = ".localization-config"
$ov.Range("B5").Value = "Not to be localized"
$ov.Range("C5").Value = "Not to be localized"

# Update the hyperlink display text to follow the new row order; the
# hyperlink Range/Id stay anchored on the same cells, only the text the
# user sees (and the r:id it keeps) needs to line up with the new values.
$ovLinks = @()
foreach ($hl in $ov.Hyperlinks) { $ovLinks += $hl }
$ovLinks[0].TextToDisplay = "b7d44263-9138-4e17-8083-9580d31ed741.md"
$ovLinks[1].TextToDisplay = "da5e8a2c-8d37-416c-94e5-65f807b540ad.md"
$ovLinks[2].TextToDisplay = "0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md"
$ovLinks[3].TextToDisplay = ".localization-config"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "b7d44263-9138-4e17-8083-9580d31ed741.md"
$zh.Range("B2").Value = "In Translation"
$zh.Range("C2").Value = "b7d44263-9138-4e17-8083-9580d31ed741.570d16d3f01576b9f1929dffc2b758d0d8b91ac2.zh-cn.xlf"
$zh.Range("D2").Value = "2016-03-02 09:23:41"
$zh.Range("G2").Value = "0001-01-01 00:00:00"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "da5e8a2c-8d37-416c-94e5-65f807b540ad.md"
$zh.Range("B3").Value = "In Translation"
$zh.Range("C3").Value = "da5e8a2c-8d37-416c-94e5-65f807b540ad.09d3fd63c6b621b81bd9fbfbe8cf3948fa5fb65c.zh-cn.xlf"
$zh.Range("D3").Value = "2016-03-02 09:23:41"
$zh.Range("G3").Value = "0001-01-01 00:00:00"
$zh.Range("H3").Value = "Include"

$zh.Range("A4").Value = "0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md"
$zh.Range("B4").Value = "Ready for handoff"
$zh.Range("C4").Value = "0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.a1ce4849b9565447c4a185bf00b08dbaf7048def.zh-cn.xlf"
$zh.Range("D4").Value = "2016-03-02 09:27:19"
$zh.Range("G4").Value = "0001-01-01 00:00:00"
$zh.Range("H4").Value = "Include"

$zh.Range("A5").Value = ".localization-config"
$zh.Range("B5").Value = "Not to be localized"
$zh.Range("D5").Value = "0001-01-01 00:00:00"
$zh.Range("G5").Value = "0001-01-01 00:00:00"
$zh.Range("H5").Value = "Ignored"

$zhLinks = @()
foreach ($hl in $zh.Hyperlinks) { $zhLinks += $hl }
$zhLinks[0].TextToDisplay = "b7d44263-9138-4e17-8083-9580d31ed741.md"
$zhLinks[1].TextToDisplay = "b7d44263-9138-4e17-8083-9580d31ed741.570d16d3f01576b9f1929dffc2b758d0d8b91ac2.zh-cn.xlf"
$zhLinks[2].TextToDisplay = "da5e8a2c-8d37-416c-94e5-65f807b540ad.md"
$zhLinks[3].TextToDisplay = "da5e8a2c-8d37-416c-94e5-65f807b540ad.09d3fd63c6b621b81bd9fbfbe8cf3948fa5fb65c.zh-cn.xlf"
$zhLinks[4].TextToDisplay = "0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md"
$zhLinks[5].TextToDisplay = "0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.a1ce4849b9565447c4a185bf00b08dbaf7048def.zh-cn.xlf"
$zhLinks[6].TextToDisplay = ".localization-config"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "b7d44263-9138-4e17-8083-9580d31ed741.md"
$de.Range("B2").Value = "In Translation"
$de.Range("C2").Value = "b7d44263-9138-4e17-8083-9580d31ed741.570d16d3f01576b9f1929dffc2b758d0d8b91ac2.de-de.xlf"
$de.Range("D2").Value = "2016-03-02 09:23:57"
$de.Range("G2").Value = "0001-01-01 00:00:00"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "da5e8a2c-8d37-416c-94e5-65f807b540ad.md"
$de.Range("B3").Value = "In Translation"
$de.Range("C3").Value = "da5e8a2c-8d37-416c-94e5-65f807b540ad.09d3fd63c6b621b81bd9fbfbe8cf3948fa5fb65c.de-de.xlf"
$de.Range("D3").Value = "2016-03-02 09:23:57"
$de.Range("G3").Value = "0001-01-01 00:00:00"
$de.Range("H3").Value = "Include"

$de.Range("A4").Value = "0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md"
$de.Range("B4").Value = "Ready for handoff"
$de.Range("C4").Value = "0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.a1ce4849b9565447c4a185bf00b08dbaf7048def.de-de.xlf"
$de.Range("D4").Value = "2016-03-02 09:27:30"
$de.Range("G4").Value = "0001-01-01 00:00:00"
$de.Range("H4").Value = "Include"

$de.Range("A5").Value = ".localization-config"
$de.Range("B5").Value = "Not to be localized"
$de.Range("D5").Value = "0001-01-01 00:00:00"
$de.Range("G5").Value = "0001-01-01 00:00:00"
$de.Range("H5").Value = "Ignored"

$deLinks = @()
foreach ($hl in $de.Hyperlinks) { $deLinks += $hl }
$deLinks[0].TextToDisplay = "b7d44263-9138-4e17-8083-9580d31ed741.md"
$deLinks[1].TextToDisplay = "b7d44263-9138-4e17-8083-9580d31ed741.570d16d3f01576b9f1929dffc2b758d0d8b91ac2.de-de.xlf"
$deLinks[2].TextToDisplay = "da5e8a2c-8d37-416c-94e5-65f807b540ad.md"
$deLinks[3].TextToDisplay = "da5e8a2c-8d37-416c-94e5-65f807b540ad.09d3fd63c6b621b81bd9fbfbe8cf3948fa5fb65c.de-de.xlf"
$deLinks[4].TextToDisplay = "0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.md"
$deLinks[5].TextToDisplay = "0bcd2ce7-5fcc-48b5-94a0-47d35f7452df.a1ce4849b9565447c4a185bf00b08dbaf7048def.de-de.xlf"
$deLinks[6].TextToDisplay = ".localization-config"

Write-Host "Report regenerated for handoff."
